# Applies scheduled runner updates to pricing/profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 925.6316
$ws.Range("I28").Value = 706.4666999999999
$ws.Range("J28").Value = 1747.5
$ws.Range("K28").Value = 706.4666999999999
$ws.Range("L28").Value = 1747.5
$ws.Range("M28").Value = -221.4666999999999
$ws.Range("N28").Value = -2717.5

$ws.Range("H111").Value = 414
$ws.Range("I111").Value = 414.83334
$ws.Range("J111").Value = 409
$ws.Range("K111").Value = 1244.50002
$ws.Range("L111").Value = 1227
$ws.Range("M111").Value = 1822.49998
$ws.Range("N111").Value = -7361

$ws.Range("H115").Value = 1600
$ws.Range("I115").Value = 1250
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 3750
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -2183
$ws.Range("N115").Value = -12134

$ws.Range("H116").Value = 131836.88
$ws.Range("I116").Value = 174532.5
$ws.Range("J116").Value = 3750
$ws.Range("K116").Value = 174532.5
$ws.Range("L116").Value = 3750
$ws.Range("M116").Value = -171090.5
$ws.Range("N116").Value = -10634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1369.625
$ws.Range("I45").Value = 1166.9375
$ws.Range("K45").Value = 1166.9375
$ws.Range("M45").Value = -789.9375

$ws.Range("H61").Value = 3115.4792
$ws.Range("I61").Value = 2220.3794
$ws.Range("J61").Value = 4481.684
$ws.Range("K61").Value = 2220.3794
$ws.Range("L61").Value = 4481.684
$ws.Range("M61").Value = -2008.3794
$ws.Range("N61").Value = -4905.684

$ws.Range("H136").Value = 3115.4792
$ws.Range("I136").Value = 2220.3794
$ws.Range("J136").Value = 4481.684
$ws.Range("K136").Value = 6661.138199999999
$ws.Range("L136").Value = 13445.052
$ws.Range("M136").Value = -4111.138199999999
$ws.Range("N136").Value = -18545.052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7907.1763
$ws.Range("I86").Value = 6701.2
$ws.Range("J86").Value = 9630
$ws.Range("K86").Value = 6701.2
$ws.Range("L86").Value = 9630
$ws.Range("M86").Value = -5578.2
$ws.Range("N86").Value = -11876

$ws.Range("H89").Value = 7907.1763
$ws.Range("I89").Value = 6701.2
$ws.Range("J89").Value = 9630
$ws.Range("K89").Value = 33506
$ws.Range("L89").Value = 48150
$ws.Range("M89").Value = -27890
$ws.Range("N89").Value = -59382

$ws.Range("H99").Value = 2940.348
$ws.Range("I99").Value = 3125.1428
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 3125.1428
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -1627.1428
$ws.Range("N99").Value = -3996

$ws.Range("H107").Value = 1806.9143
$ws.Range("I107").Value = 1562.2
$ws.Range("J107").Value = 3275.2
$ws.Range("K107").Value = 1562.2
$ws.Range("L107").Value = 3275.2
$ws.Range("M107").Value = 357.8
$ws.Range("N107").Value = -7115.2

$ws.Range("H134").Value = 22096.691
$ws.Range("I134").Value = 29646.432
$ws.Range("J134").Value = 6577.778
$ws.Range("K134").Value = 88939.296
$ws.Range("L134").Value = 19733.334
$ws.Range("M134").Value = -86404.296
$ws.Range("N134").Value = -24803.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1392.8572
$ws.Range("I16").Value = 1150
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1150
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -863
$ws.Range("N16").Value = -2574

$ws.Range("H58").Value = 1449.4878
$ws.Range("I58").Value = 922.2222
$ws.Range("K58").Value = 922.2222
$ws.Range("M58").Value = -719.2222

$ws.Range("H113").Value = 1392.8572
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 2238.25
$ws.Range("I132").Value = 1200.2667
$ws.Range("J132").Value = 3154.1177
$ws.Range("K132").Value = 3600.800099999999
$ws.Range("L132").Value = 9462.3531
$ws.Range("M132").Value = -1070.800099999999
$ws.Range("N132").Value = -14522.3531

$ws.Range("H136").Value = 1449.4878
$ws.Range("I136").Value = 922.2222
$ws.Range("K136").Value = 2766.6666
$ws.Range("M136").Value = -216.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3061
$ws.Range("I3").Value = 2733.25
$ws.Range("J3").Value = 3716.5
$ws.Range("K3").Value = 8199.75
$ws.Range("L3").Value = 11149.5
$ws.Range("M3").Value = -8087.75
$ws.Range("N3").Value = -11373.5

$ws.Range("H68").Value = 427.27274
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 300
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 900
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -2522

$ws.Range("H71").Value = 427.27274
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 2700
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -10812

$ws.Range("H129").Value = 1583.3334
$ws.Range("I129").Value = 596.9
$ws.Range("J129").Value = 2816.375
$ws.Range("K129").Value = 1790.7
$ws.Range("L129").Value = 8449.125
$ws.Range("M129").Value = 3209.3
$ws.Range("N129").Value = -18449.125

$ws.Range("H131").Value = 882.2787
$ws.Range("J131").Value = 910.6842
$ws.Range("L131").Value = 2732.0526
$ws.Range("N131").Value = -12812.0526

$ws.Range("H136").Value = 4352.3
$ws.Range("I136").Value = 4165
$ws.Range("J136").Value = 4633.25
$ws.Range("K136").Value = 12495
$ws.Range("L136").Value = 13899.75
$ws.Range("M136").Value = -7395
$ws.Range("N136").Value = -24099.75

$ws.Range("H138").Value = 2200.0557
$ws.Range("I138").Value = 1897.2727
$ws.Range("J138").Value = 2675.8572
$ws.Range("K138").Value = 5691.8181
$ws.Range("L138").Value = 8027.571599999999
$ws.Range("M138").Value = -551.8181000000004
$ws.Range("N138").Value = -18307.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 11492.667
$ws.Range("I107").Value = 12841.75
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 12841.75
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = -10921.75
$ws.Range("N107").Value = -4540

$ws.Range("H132").Value = 5729.478
$ws.Range("I132").Value = 7730
$ws.Range("J132").Value = 4190.615
$ws.Range("K132").Value = 23190
$ws.Range("L132").Value = 12571.845
$ws.Range("M132").Value = -20660
$ws.Range("N132").Value = -17631.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6039.375
$ws.Range("I136").Value = 3536.875
$ws.Range("J136").Value = 8541.875
$ws.Range("K136").Value = 10610.625
$ws.Range("L136").Value = 25625.625
$ws.Range("M136").Value = -8060.625
$ws.Range("N136").Value = -30725.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 28371028
$ws.Range("I136").Value = 41668550
$ws.Range("J136").Value = 14495349
$ws.Range("K136").Value = 125005650
$ws.Range("L136").Value = 43486047
$ws.Range("M136").Value = -125003100
$ws.Range("N136").Value = -43491147
